# Automatische test-sync: 2025-06-19 14:00:10
# Adds the new mail-log entry for 2025-06-19 13:58:11 to the "Logs" sheet
# and refreshes the dependent "Dashboard" summary/conditional formatting.

$wb = $excel.ActiveWorkbook

# --- 1. Append the new log row (row 13) on the "Logs" sheet -----------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A13").Value = "Klacht over levering"
$logs.Range("B13").Value = "mailmind.test@zohomail.eu"
$logs.Range("C13").Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$logs.Range("D13").Value = "Klacht"
$logs.Range("F13").Value = "2025-06-19 13:58:11"
$logs.Range("G13").Value = "Nee"

# --- 2. Extend the conditional formatting ranges so they include row 13 -----
$fcCategorie = $logs.Range("D2:D12").FormatConditions.Item(1)
$fcCategorie.ModifyAppliesToRange($logs.Range("D2:D13"))

$fcBeantwoord = $logs.Range("G2:G12").FormatConditions.Item(1)
$fcBeantwoord.ModifyAppliesToRange($logs.Range("G2:G13"))

# --- 3. Update the "Dashboard" summary count for the "Klacht" category ------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B4").Value = 2
